$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '59.494.67'
$ws.Range('E2').Value = '  -5.44%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.444.05'
$ws.Range('E3').Value = '  -8.79%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '535.35'
$ws.Range('E5').Value = '  -2.85%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '147.14'
$ws.Range('E6').Value = '  -6.38%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.568'
$ws.Range('E8').Value = '  -3.32%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0991'
$ws.Range('E9').Value = '  -5.84%  '
$ws.Range('E10').Value = '  -2.20%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.48'
$ws.Range('E11').Value = '  +7.63%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.350'
$ws.Range('E12').Value = '  -4.46%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.882.86'
$ws.Range('E13').Value = '  -8.61%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '24.11'
$ws.Range('E14').Value = '  -7.06%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '59.447.08'
$ws.Range('E15').Value = '  -5.33%  '
$ws.Range('E16').Value = '  -5.87%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.491.73'
$ws.Range('E17').Value = '  -7.13%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.18'
$ws.Range('E18').Value = '  -5.78%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.38'
$ws.Range('E19').Value = '  -3.91%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '325.04'
$ws.Range('E20').Value = '  -5.04%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.971'
$ws.Range('E21').Value = '  -2.68%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.71'
$ws.Range('E22').Value = '  -9.51%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.467'
$ws.Range('E23').Value = '  -7.08%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '60.47'
$ws.Range('E24').Value = '  -4.49%  '
$ws.Range('E25').Value = '  -3.53%  '
$ws.Range('E26').Value = '  -2.44%  '
$ws.Range('E27').Value = '  -4.29%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.30'
$ws.Range('E28').Value = '  -2.56%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.88'
$ws.Range('E29').Value = '  -1.87%  '
$ws.Range('E30').Value = '  -5.63%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0₃0766'
$ws.Range('E31').Value = '  -10.14%  '
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '156.09'
$ws.Range('E33').Value = '  -6.42%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.54'
$ws.Range('E34').Value = '  -5.04%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '18.36'
$ws.Range('E35').Value = '  -5.83%  '
$ws.Range('E36').Value = '  -5.18%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '315.67'
$ws.Range('E38').Value = '  -6.46%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.79'
$ws.Range('E39').Value = '  -6.26%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.852'
$ws.Range('E40').Value = '  -8.26%  '
$ws.Range('E41').Value = '  -3.57%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.73'
$ws.Range('E42').Value = '  -5.14%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.997'
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('E44').Value = '  -2.85%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.582'
$ws.Range('E45').Value = '  -5.38%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0939'
$ws.Range('E46').Value = '  -3.26%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0525'
$ws.Range('E47').Value = '  -6.19%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '18.59'
$ws.Range('E48').Value = '  -7.91%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0230'
$ws.Range('E49').Value = '  -3.86%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '19.04'
$ws.Range('E50').Value = '  -7.79%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.991.67'
$ws.Range('E51').Value = '  -4.48%  '
